# Daily attendance processing - 2025-10-14 14:23:37
# Normalize "Recorded By" (column G) ordering: move "System" after the
# email address for the rows that currently read "System, <email>".

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Session Analysis Results")

$rows = @(3,6,10,11,12,13,14,15,17,18,19,30,33,37,38,39,40,41,42,44,45,46,57,60,64,65,66,67,68,69,71,72,73,86,87,88,89,93,95,96,97,112,113,114,115,119,121,122,123,138,139,140,141,145,147,148,149)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
